$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of row 19 and row 20 (columns A and B)
$a19 = $ws.Range("A19").Value()
$b19 = $ws.Range("B19").Value()
$a20 = $ws.Range("A20").Value()
$b20 = $ws.Range("B20").Value()

$ws.Range("A19").Value = $a20
$ws.Range("B19").Value = $b20
$ws.Range("A20").Value = $a19
$ws.Range("B20").Value = $b19

# Update the active cell selection to A26
$ws.Range("A26").Select()

# Update the workbook window position/size
$wb.Windows.Item(1).Left = -105
$wb.Windows.Item(1).Top = 0
$wb.Windows.Item(1).Width = 10455
$wb.Windows.Item(1).Height = 10905
